$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 1689
$ws1.Range("F8").Value = 1894
$ws1.Range("F11").Value = 774
$ws1.Range("F13").Value = 178
$ws1.Range("F19").Value = 6862
$ws1.Range("F25").Value = 321

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 1689
$ws4.Range("F11").Value = 1894
$ws4.Range("F14").Value = 774
$ws4.Range("F17").Value = 178
$ws4.Range("F22").Value = 6862
$ws4.Range("F29").Value = 321
